# Update "F2", "F5", "F6", "F8" (想去人数 / headcount) figures on both the
# "展览" sheet and the "全部类型" sheet, which carry duplicate rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 249
    $ws.Range("F5").Value = 6506
    $ws.Range("F6").Value = 5268
    $ws.Range("F8").Value = 64
}
